$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text values that can look numeric (e.g. "1.001", "326.55").
# Force text storage (matching the source t="inlineStr" cells) by switching
# NumberFormat to "@" before assigning the value, then restore the cells
# original Style afterwards so no stray style index lingers on the cell.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.192.50"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  +3.50%  "
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.03"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  +0.26%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.23%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.55"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +3.58%  "
$ws.Range("E6").Value = "  -0.18%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +0.75%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4014"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  +0.55%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.66"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -0.06%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.24"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +13.35%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.441"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +3.12%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.908.86"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +0.47%  "
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.354"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  -0.25%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.80"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.85%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001113"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +0.47%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06716"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("E21").Value = "  -0.09%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.999"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.24%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.196.12"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +3.47%  "
$ws.Range("E24").Value = "  +0.78%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.207"
$ws.Range("D25").Style = $origStyle
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.124.46"
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.68"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +3.48%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.63"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +1.40%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.385"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -1.80%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.31"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +1.81%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.098"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("E32").Value = "  +1.44%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.048"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +2.24%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.701"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.95%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02495"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  -0.56%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2199"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +0.41%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.195"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  +0.07%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.91"
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.810"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -2.90%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6509"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("E43").Value = "  -0.07%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6118"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +1.47%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +0.08%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.715"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +1.17%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.060"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  +1.23%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.06"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  -0.31%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.19"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +2.06%  "
